{"js": "// Find the \"BranchMilkLog (...)\" definition paragraph and apply a single\n// underline to the \"IsMorning\" and \"LocalDID\" field names within it\n// (matching the author's edit for the BranchMilkLog table update).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"BranchMilkLog (\") === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  // Underline \"IsMorning\" inside this paragraph only.\n  const isMorningRanges = target.search(\"IsMorning\", { matchCase: true });\n  isMorningRanges.load(\"text\");\n  await context.sync();\n  if (isMorningRanges.items.length > 0) {\n    isMorningRanges.items[0].font.underline = \"Single\";\n  }\n\n  // Underline \"LocalDID\" inside this paragraph only (keeps its existing\n  // grey shading, just adds the underline on top of it).\n  const localDidRanges = target.search(\"LocalDID\", { matchCase: true });\n  localDidRanges.load(\"text\");\n  await context.sync();\n  if (localDidRanges.items.length > 0) {\n    localDidRanges.items[0].font.underline = \"Single\";\n  }\n\n  await context.sync();\n}\n", "ps1": "# Find the \"BranchMilkLog (...)\" definition paragraph and apply a single\n# underline to the \"IsMorning\" and \"LocalDID\" field names within it\n# (matching the author's edit for the BranchMilkLog table update).\n\n$d = $word.ActiveDocument\n$wdUnderlineSingle = 1\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"BranchMilkLog (\")) {\n\n        # Underline \"IsMorning\" inside this paragraph only.\n        $rngMorning = $p.Range\n        $rngMorning.Find.ClearFormatting()\n        $rngMorning.Find.Text = \"IsMorning\"\n        $rngMorning.Find.MatchCase = $true\n        $rngMorning.Find.Execute() | Out-Null\n        if ($rngMorning.Find.Found) {\n            $rngMorning.Font.Underline = $wdUnderlineSingle\n        }\n\n        # Underline \"LocalDID\" inside this paragraph only (keeps its\n        # existing grey shading, just adds the underline on top of it).\n        $rngLocal = $p.Range\n        $rngLocal.Find.ClearFormatting()\n        $rngLocal.Find.Text = \"LocalDID\"\n        $rngLocal.Find.MatchCase = $true\n        $rngLocal.Find.Execute() | Out-Null\n        if ($rngLocal.Find.Found) {\n            $rngLocal.Font.Underline = $wdUnderlineSingle\n        }\n\n        break\n    }\n}\n"}
